# Data update using git
# Update "Pagos" (col F) and "Inscrições homologadas" (col H) values for
# four rows in the "Inscricoes" worksheet/table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 17: Pagos 87 -> 88, Homologadas 119 -> 120
$ws.Range("F17").Value = 88
$ws.Range("H17").Value = 120

# Row 50: Pagos 14 -> 15, Homologadas 23 -> 24
$ws.Range("F50").Value = 15
$ws.Range("H50").Value = 24

# Row 72: Pagos 35 -> 36, Homologadas 46 -> 47
$ws.Range("F72").Value = 36
$ws.Range("H72").Value = 47

# Row 83: Pagos 4 -> 5, Homologadas 11 -> 12
$ws.Range("F83").Value = 5
$ws.Range("H83").Value = 12
